$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Add new transaction row 37 (copy formatting down from the last data row, then set values) ---
$ws1.Range("A36:F36").Copy()
$ws1.Range("A37:F37").PasteSpecial(-4122)

$ws1.Range("A37").Value = "MI"
$ws1.Range("B37").Value = "LCWD"
$ws1.Range("C37").Value = 45231
$ws1.Range("D37").Value = 46
$ws1.Range("E37").Value = 13.638
$ws1.Range("F37").Value = 1.5

# --- Add 5 blank (but formatted) rows below it, 38-42 ---
$ws1.Range("A37:B37").Copy()
$ws1.Range("A38:C42").PasteSpecial(-4122)
$ws1.Range("D37").Copy()
$ws1.Range("D38:D42").PasteSpecial(-4122)
$ws1.Range("E37").Copy()
$ws1.Range("E38:E42").PasteSpecial(-4122)

# --- Resize the table / autofilter to include the new data row ---
$lo = $ws1.ListObjects.Item(1)
$lo.Resize($ws1.Range("A1:F37"))

# --- Update selections / active cells on both sheets ---
# (select sheet2's cell first so sheet1 ends up as the active/selected tab)
$null = $ws2.Range("C15").Select()
$null = $ws1.Range("C25").Select()
